$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.369.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "'1.567.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'290.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'0.3754"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").Value = "'49.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'0.3395"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'0.07556"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'1.132"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'20.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "'5.957"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "'6.915"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'1.565.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.00001128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'89.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "'0.06757"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'16.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'6.182"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "'11.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "'22.373.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.379"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'2.695"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'147.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'5.023"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'125.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'1.741.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'2.022"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "'6.038"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "'0.9844"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "'9.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'1.417"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.11%  "
$ws.Range("D37").Value = "'0.08443"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'0.02486"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").Value = "'0.2283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").Value = "'0.06465"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'5.391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").Value = "'0.6288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").Value = "'11.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'14.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'3.799"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'0.5907"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'2.064"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "'1.262"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'124.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "'0.07329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.78%  "
